# Auto-committed on 2023/05/19 週五 17:22:58.97
# Updates a batch of "最後修改時間" (last-modified) timestamps for several
# GenTable entries, and adds a brand-new row for the "TxControl" table
# (inserted alphabetically between TxBizDate and TxCruiser, in the
# "XX-系統" category block).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Refresh the "last modified" timestamp string for a batch of tables.
#    Each (old timestamp -> new timestamp) pair below is unique within
#    the sheet, so Cells.Find locates the single matching cell reliably.
# ---------------------------------------------------------------------
$timestampUpdates = @(
    @("2023年05月04日 16:53:05", "2023年05月18日 10:42:39"),  # LoanBorMain
    @("2023年05月02日 10:58:51", "2023年05月19日 12:19:18"),  # AchAuthLog
    @("2023年02月15日 20:06:36", "2023年05月19日 13:24:23"),  # InsuRenew
    @("2023年04月28日 09:31:55", "2023年05月19日 12:19:14"),  # PostAuthLog
    @("2022年01月20日 11:27:22", "2023年05月19日 14:58:04"),  # PfCoOfficer
    @("2021年09月27日 13:48:53", "2023年05月19日 16:04:06"),  # PfCoOfficerLog
    @("2022年01月20日 11:28:43", "2023年05月18日 18:13:01"),  # PfIntranetAdjust
    @("2022年01月20日 11:30:02", "2023年05月18日 12:44:47"),  # PfReward
    @("2021年12月02日 16:50:40", "2023年05月18日 17:28:31"),  # PfRewardMedia
    @("2023年01月17日 10:50:54", "2023年05月15日 16:28:26"),  # CdCode
    @("2022年10月12日 16:34:36", "2023年05月15日 17:25:02"),  # JobMain
    @("2020年07月13日 09:41:39", "2021年07月15日 10:15:27")   # TxHoliday
)

foreach ($pair in $timestampUpdates) {
    $oldValue = $pair[0]
    $newValue = $pair[1]
    $cell = $ws.Cells.Find($oldValue)
    $cell.Value = $newValue
}

# ---------------------------------------------------------------------
# 2) Insert a new row for the "TxControl" table ("作業流程控制檔"),
#    right before the existing "TxCruiser" row (keeps the XX-系統 block
#    alphabetically ordered), and fill in its five columns.
# ---------------------------------------------------------------------
$cruiserCell = $ws.Cells.Find("TxCruiser")
$cruiserRow = $cruiserCell.Row

$ws.Rows($cruiserRow).Insert()

$newRow = $cruiserRow
$category = $ws.Cells.Item($newRow + 1, 1).Value2

$ws.Cells.Item($newRow, 1).Value = $category
$ws.Cells.Item($newRow, 2).Value = "TxControl"
$ws.Cells.Item($newRow, 3).Value = "作業流程控制檔"
$ws.Cells.Item($newRow, 4).Formula = '=HYPERLINK("[\\192.168.10.16\St1Share(NAS)\SKL\DB\GenTables\XX-系統\TxControl.xlsx]DBD!A1", "連結")'
$ws.Cells.Item($newRow, 5).Value = "2023年05月19日 12:32:26"
